# Add a new "TEST" row (row 2) under the existing header row on the active
# sheet, matching the author's uploaded revision:
#   - new shared string "TEST"
#   - row 2, columns A:G all contain "TEST"
#   - row 2 keeps the same visual style as the rest of the data columns
#     (DejaVu Sans Mono 8pt, vertical-top, wrap text - i.e. the same
#     formatting already used by columns A:E)
#   - the active cell/selection moves to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row of data.
$ws.Range("A2:G2").Value = "TEST"

# Columns F:G default to a plain numeric style (no wrap/font override), so
# bring them in line with the rest of the row's formatting explicitly.
$ws.Range("F2:G2").Font.Name = "DejaVu Sans Mono"
$ws.Range("F2:G2").Font.Size = 8
$ws.Range("F2:G2").VerticalAlignment = -4160   # xlVAlignTop
$ws.Range("F2:G2").WrapText = $true

# Move the selection/active cell to A2, as in the saved workbook.
$ws.Range("A2").Select() | Out-Null
